# "app and report work"
#
# - Week 5: shift the three trailing dates forward by 4 weeks (28 days) and
#   move the selection down to B11.
# - week 6: add four new progress rows (8-11), the last one bold, and move
#   the selection/dimension accordingly.
# - Add a new "Week 7" sheet at the end of the tab strip with a single new
#   dated entry, and make it the active sheet/selection.

$wb = $excel.ActiveWorkbook

# --- "Week 5" sheet: rework dates for rows 8-10, move selection ---
$ws5 = $wb.Worksheets.Item("Week 5")
$ws5.Range("B8").Value = 42797
$ws5.Range("B9").Value = 42798
$ws5.Range("B10").Value = 42799
$ws5.Range("B11").Select() | Out-Null

# --- "week 6" sheet: append new progress rows 8-11 ---
$ws6 = $wb.Worksheets.Item("week 6")

$ws6.Range("B8").Value = 42775
$ws6.Range("B8").NumberFormat = "d-mmm"
$ws6.Range("C8").Value = "data about location pulled into web app"

$ws6.Range("C9").Value = "Graph added to app using android graph view"

$ws6.Range("B10").Value = 42776
$ws6.Range("B10").NumberFormat = "d-mmm"
$ws6.Range("C10").Value = "Location data plotted on graph"

$ws6.Range("C11").Value = "Test pulling data and displaying graph"
$ws6.Range("C11").Font.Bold = $true

$ws6.Range("C11").Select() | Out-Null

# --- New "Week 7" sheet, inserted after "week 6" (end of tab strip) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws7.Name = "Week 7"

$ws7.Range("B3").Value = 42779
$ws7.Range("B3").NumberFormat = "d-mmm"
$ws7.Range("C3").Value = "Work on report"

$ws7.Range("C3").Select() | Out-Null
